# Applies the OOXML-level diff:
#  - merges a few runs that were needlessly split
#  - fixes the "woren" -> "worden" typo
#  - wraps several technical / mistyped terms with <w:proofErr> spell-check
#    markers (splitting the surrounding run, exactly like Word's own
#    on-the-fly proofer would have done), leaving paragraph/run formatting
#    (pPr, rPr) untouched
#  - tidies up a run of duplicated whitespace runs around the _GoBack
#    bookmark in the final paragraph

$d = $word.ActiveDocument

function Set-ParaInnerXml {
    param($doc, [int]$paraIndex, [string]$innerXml)

    $para = $doc.Paragraphs.Item($paraIndex)
    $full = $para.Range
    # Exclude the trailing paragraph mark so the paragraph's own <w:pPr>
    # (list numbering, spacing, rPr-of-mark, ...) is left completely alone;
    # InsertXML only replaces the content of the range it is called on.
    $body = $doc.Range($full.Start, $full.End - 1)

    $pkg = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>'
    $pkg += '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">'
    $pkg += '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">'
    $pkg += '<pkg:xmlData>'
    $pkg += '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">'
    $pkg += '<w:body><w:p>'
    $pkg += $innerXml
    $pkg += '</w:p></w:body></w:document>'
    $pkg += '</pkg:xmlData></pkg:part></pkg:package>'

    $body.InsertXML($pkg)
}

# --- 1. Title: merge "Bijlage Onderzoeksrapport " + "Applicatie" ----------
$xml1 = '<w:r><w:t>Bijlage Onderzoeksrapport Applicatie</w:t></w:r>'
Set-ParaInnerXml $d 1 $xml1

# --- 2. "Jullie" + " hebben de applicatie ... en meer." --------------------
$xml2 = '<w:r><w:t>Jullie hebben de applicatie ontvangen waar jullie mee gaan werken. '
$xml2 += 'Dat betekent dat deze volledig wordt uitgeplozen voordat jullie hier mee aan de slag gaan. '
$xml2 += 'Zowel de applicatie als de database wordt door jullie onder de loep genomen. '
$xml2 += 'Op basis van dit onderzoek gaan jullie een adviesvoorstel doen om deze applicatie te optimaliseren. '
$xml2 += 'Denk bijvoorbeeld aan het toevoegen van de finalewedstrijden, bijhouden van de scores per speler en meer.</w:t></w:r>'
Set-ParaInnerXml $d 3 $xml2

# --- 3. "In het onderzoeksrapport" + " " + "worden de volgende ..." -------
$xml3 = '<w:r><w:t>In het onderzoeksrapport worden de volgende leervragen beantwoord:</w:t></w:r>'
Set-ParaInnerXml $d 5 $xml3

# --- 4. Leervraag 1: flag selectWedstrijden.php ----------------------------
$xml4 = '<w:r><w:t xml:space="preserve">Beschrijf hoe de resetknop werkt in </w:t></w:r>'
$xml4 += '<w:proofErr w:type="spellStart"/>'
$xml4 += '<w:r><w:t>selectWedstrijden.php</w:t></w:r>'
$xml4 += '<w:proofErr w:type="spellEnd"/>'
Set-ParaInnerXml $d 6 $xml4

# --- 5. Leervraag 2: flag "stream", fix woren -> worden --------------------
$xml5 = '<w:r><w:t xml:space="preserve">In de live </w:t></w:r>'
$xml5 += '<w:proofErr w:type="spellStart"/>'
$xml5 += '<w:r><w:t>stream</w:t></w:r>'
$xml5 += '<w:proofErr w:type="spellEnd"/>'
$xml5 += '<w:r><w:t xml:space="preserve"> pagina wor</w:t></w:r>'
$xml5 += '<w:r><w:t>d</w:t></w:r>'
$xml5 += '<w:r><w:t>en de statistieken real-time bijgewerkt. Waarom is daarvoor gekozen?</w:t></w:r>'
Set-ParaInnerXml $d 7 $xml5

# Shared rPr block used by the "Lijstalinea" leervragen further down.
$rPr = '<w:rPr><w:rFonts w:eastAsia="Times New Roman" w:cs="Helvetica"/><w:color w:val="000000"/><w:sz w:val="21"/><w:szCs w:val="21"/><w:lang w:eastAsia="nl-NL"/></w:rPr>'

# --- 6. "Wat is de query ... gespeelde wedstrjiden te selecteren?" --------
$xml6 = '<w:r>'
$xml6 += $rPr
$xml6 += '<w:t xml:space="preserve">Wat is de query die gebruikt wordt om een lijst van gespeelde </w:t></w:r>'
$xml6 += '<w:proofErr w:type="spellStart"/>'
$xml6 += '<w:r>'
$xml6 += $rPr
$xml6 += '<w:t>wedstrjiden</w:t></w:r>'
$xml6 += '<w:proofErr w:type="spellEnd"/>'
$xml6 += '<w:r>'
$xml6 += $rPr
$xml6 += '<w:t xml:space="preserve"> te selecteren?</w:t></w:r>'
Set-ParaInnerXml $d 9 $xml6

# --- 7. "In welke condities ... addScore.php behandeld?" ------------------
$xml7 = '<w:r>'
$xml7 += $rPr
$xml7 += '<w:t xml:space="preserve">In welke condities wordt momenteel de score van </w:t></w:r>'
$xml7 += '<w:proofErr w:type="spellStart"/>'
$xml7 += '<w:r>'
$xml7 += $rPr
$xml7 += '<w:t>addScore.php</w:t></w:r>'
$xml7 += '<w:proofErr w:type="spellEnd"/>'
$xml7 += '<w:r>'
$xml7 += $rPr
$xml7 += '<w:t xml:space="preserve"> behandeld?</w:t></w:r>'
Set-ParaInnerXml $d 10 $xml7

# --- 8. "Hoe weet de applicatie ... selectWedstrijd.php?" -----------------
$xml8 = '<w:r>'
$xml8 += $rPr
$xml8 += '<w:t xml:space="preserve">Hoe weet de applicatie welke wedstrijd ik heb gekozen in </w:t></w:r>'
$xml8 += '<w:proofErr w:type="spellStart"/>'
$xml8 += '<w:r>'
$xml8 += $rPr
$xml8 += '<w:t>selectWedstrijd.php</w:t></w:r>'
$xml8 += '<w:proofErr w:type="spellEnd"/>'
$xml8 += '<w:r>'
$xml8 += $rPr
$xml8 += '<w:t>?</w:t></w:r>'
Set-ParaInnerXml $d 11 $xml8

# --- 9. "...n de streamview zie je: 'nu speelt'..." (2nd run only) --------
$xml9 = '<w:r><w:rPr><w:rFonts w:eastAsia="Times New Roman" w:cs="Helvetica"/><w:color w:val="000000"/><w:sz w:val="21"/><w:szCs w:val="21"/><w:lang w:eastAsia="nl-NL"/></w:rPr><w:t>I</w:t></w:r>'
$xml9 += '<w:r>'
$xml9 += $rPr
$xml9 += '<w:t xml:space="preserve">n de </w:t></w:r>'
$xml9 += '<w:proofErr w:type="spellStart"/>'
$xml9 += '<w:r>'
$xml9 += $rPr
$xml9 += '<w:t>streamview</w:t></w:r>'
$xml9 += '<w:proofErr w:type="spellEnd"/>'
$xml9 += '<w:r>'
$xml9 += $rPr
$xml9 += '<w:t xml:space="preserve">'
$xml9 += " zie je: 'nu speelt'... hoe komt deze data tot stand?"
$xml9 += '</w:t></w:r>'
Set-ParaInnerXml $d 12 $xml9

# --- 10. "Leg in eigen woorden uit wat de 'Join' statement ..." -----------
$xml10 = '<w:r>'
$xml10 += $rPr
$xml10 += '<w:t>Leg in eigen woorden uit wat de '
$xml10 += [char]0x2018
$xml10 += '</w:t></w:r>'
$xml10 += '<w:proofErr w:type="spellStart"/>'
$xml10 += '<w:r>'
$xml10 += $rPr
$xml10 += '<w:t>Join</w:t></w:r>'
$xml10 += '<w:proofErr w:type="spellEnd"/>'
$xml10 += '<w:r>'
$xml10 += $rPr
$xml10 += '<w:t>'
$xml10 += [char]0x2019
$xml10 += ' statement betekent in SQL.</w:t></w:r>'
Set-ParaInnerXml $d 14 $xml10

# --- 11. "Hoe wordt de streaming weergegeven in livestream.html?" ---------
$xml11 = '<w:r>'
$xml11 += $rPr
$xml11 += '<w:t xml:space="preserve">Hoe wordt de </w:t></w:r>'
$xml11 += '<w:proofErr w:type="spellStart"/>'
$xml11 += '<w:r>'
$xml11 += $rPr
$xml11 += '<w:t>streaming</w:t></w:r>'
$xml11 += '<w:proofErr w:type="spellEnd"/>'
$xml11 += '<w:r>'
$xml11 += $rPr
$xml11 += '<w:t xml:space="preserve"> weergegeven in livestream.html?</w:t></w:r>'
Set-ParaInnerXml $d 15 $xml11

# --- 12. Final paragraph: flag "invoerforms", collapse duplicate spaces ---
$xml12 = '<w:r><w:t>Op welke manieren kan de code verbeterd worden? Op welke manier kan de folderstructuur beter worden weggezet?</w:t></w:r>'
$xml12 += '<w:r><w:t xml:space="preserve"> Hoe kunnen de </w:t></w:r>'
$xml12 += '<w:proofErr w:type="spellStart"/>'
$xml12 += '<w:r><w:t>invoerforms</w:t></w:r>'
$xml12 += '<w:proofErr w:type="spellEnd"/>'
$xml12 += '<w:r><w:t xml:space="preserve"> beter gevalideerd worden?</w:t></w:r>'
$xml12 += '<w:r><w:t xml:space="preserve"> </w:t></w:r>'
$xml12 += '<w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/>'
$xml12 += '<w:r><w:t>Kortom, al jullie voorstellen tot optimalisatie van de applicatie komt hierin tot stand.</w:t></w:r>'
Set-ParaInnerXml $d 19 $xml12
